# Insert a new weekly record as row 36 (pushing the existing rows 36-116
# down to 37-117) on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 36..116 down by one row to make room for the new record.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record's data.
$ws.Cells.Item(36, 1).Value = 5
$ws.Cells.Item(36, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value = "Maule"
$ws.Cells.Item(36, 4).Value = 45162
$ws.Cells.Item(36, 5).Value = 7
$ws.Cells.Item(36, 6).Value = 100112040
$ws.Cells.Item(36, 7).Value = "Cilantro"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 8000
$ws.Cells.Item(36, 12).Value = 8000
$ws.Cells.Item(36, 13).Value = 8000
$ws.Cells.Item(36, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 222
$ws.Cells.Item(36, 17).Value = 36
$ws.Cells.Item(36, 18).Value = "Hortaliza"
